$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PRODUCT BACKLOG")

$ws.Range("D3").Value = "Adel"
$ws.Range("D4").Value = "Adel"
$ws.Range("D6").Value = "Donovan"
$ws.Range("D8").Value = "Alihan"
$ws.Range("D11").Value = "Adel"
$ws.Range("D12").Value = "Donovan"
$ws.Range("D13").Value = "Adel"
$ws.Range("D17").Value = "Manojlo"
$ws.Range("D18").Value = "Adel"
$ws.Range("D20").Value = "Adel"
$ws.Range("D22").Value = "Alihan"

$ws.Activate()
$ws.Range("E21").Select()
